# Add a "Turkey" worksheet (Zettler Turkey template test data), based on
# the existing "Spain" sheet layout, after the last sheet ("Spain").

$wb = $excel.ActiveWorkbook
$spain = $wb.Worksheets.Item("Spain")

# Duplicate Spain (keeps styles / merged cells / column widths identical)
# right after itself, then rename the duplicate to "Turkey".
$spain.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Spain's product list has three extra "Lite" rows (8:10) and a block of
# Swiss-only MZX/MX/ZX rows (22:30 once the former are removed) that the
# Turkey sheet does not carry - drop them so the row layout matches.
$turkey.Rows("8:10").Delete()
$turkey.Rows("22:30").Delete()

# Market name / ticket reference specific to Turkey.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3290"

# Rows 3:5 on Spain carry an explicit (wrapped-text) row height; Turkey's
# sheet uses the default row height instead.
$turkey.Rows("3:5").AutoFit()

# Column D is a little wider on the Turkey sheet than on Spain's.
$turkey.Columns("D").ColumnWidth = 26

# Turkey becomes the active / selected sheet, scrolled so F18 is selected.
$turkey.Activate()
$turkey.Range("F18").Select()

# Spain itself is left scrolled further down with the whole table selected.
$spain.Range("A1:D36").Select()

$turkey.Activate()
